$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns D (Fecha), J (Volumen), K (Precio mínimo),
# L (Precio máximo), M (Precio promedio ponderado), P (Precio $/Kg)
# for rows 3 through 9 (the weekly data got re-dated / re-permuted).

$data = @{
    3 = @{ D = 44397; J = 140; K = 12500; L = 13000; M = 12750; P = 981 }
    4 = @{ D = 44320; J = 160; K = 19000; L = 20000; M = 19500; P = 1500 }
    5 = @{ D = 44379; J = 120; K = 12000; L = 13000; M = 12667; P = 974 }
    6 = @{ D = 44159; J = 100; K = 23000; L = 24000; M = 23500; P = 1808 }
    7 = @{ D = 44229; J = 120; K = 44000; L = 45000; M = 44500; P = 3423 }
    8 = @{ D = 44389; J = 120; K = 12000; L = 13000; M = 12500; P = 962 }
    9 = @{ D = 44469; J = 140; K = 13000; L = 14000; M = 13500; P = 1038 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value = $vals.D
    $ws.Cells.Item($row, 10).Value = $vals.J
    $ws.Cells.Item($row, 11).Value = $vals.K
    $ws.Cells.Item($row, 12).Value = $vals.L
    $ws.Cells.Item($row, 13).Value = $vals.M
    $ws.Cells.Item($row, 16).Value = $vals.P
}
